$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current row 2 / row 3 values (columns A and B) before overwriting
$a2 = $ws.Range("A2").Value()
$b2 = $ws.Range("B2").Value()
$a3 = $ws.Range("A3").Value()
$b3 = $ws.Range("B3").Value()

# Swap rows 2 and 3
$ws.Range("A2").Value = $a3
$ws.Range("B2").Value = $b3
$ws.Range("A3").Value = $a2
$ws.Range("B3").Value = $b2

# New column C: header + values
$ws.Range("C1").Value = "Kommune_Nr"
$ws.Range("C2").Value = 11111
$ws.Range("C3").Value = 22222

# New empty row 4 (keeps the bordered style by copying format from row 3)
$ws.Range("A3:B3").Copy() | Out-Null
$ws.Range("A4:B4").PasteSpecial(-4122) | Out-Null
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = ""

$ws.Range("C4").Select() | Out-Null
